$d = $word.ActiveDocument

$pairs = @(
    @{old="298÷7="; new="347÷6="},
    @{old="632÷8="; new="171÷7="},
    @{old="942÷6="; new="821÷4="},
    @{old="854÷6="; new="209÷3="},
    @{old="332÷7="; new="276÷3="},
    @{old="316÷2="; new="223÷7="},
    @{old="908÷3="; new="499÷7="},
    @{old="387÷6="; new="868÷4="},
    @{old="870÷5="; new="543÷5="},
    @{old="668÷2="; new="414÷3="},
    @{old="306÷6="; new="841÷3="},
    @{old="133÷8="; new="260÷7="},
    @{old="941÷8="; new="492÷7="},
    @{old="751÷3="; new="165÷7="},
    @{old="489÷5="; new="698÷6="},
    @{old="843÷5="; new="200÷5="},
    @{old="653÷4="; new="879÷4="},
    @{old="180÷4="; new="174÷2="},
    @{old="920÷9="; new="548÷3="},
    @{old="757÷6="; new="821÷2="},
    @{old="263÷6="; new="915÷3="},
    @{old="783÷2="; new="629÷9="},
    @{old="616÷6="; new="742÷9="},
    @{old="888÷6="; new="782÷7="},
    @{old="560÷4="; new="539÷5="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
